# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-25 10:58:07"
$wsZhCn.Range("G3").Value = "2016-01-25 10:58:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-25 10:58:16"
$wsDeDe.Range("G3").Value = "2016-01-25 10:59:05"
